# Fruta / hortaliza, semanal
# Update the weekly price records: the Fecha (date), Calidad, Volumen,
# Precio minimo/maximo/promedio ponderado and Precio $/Kg columns are
# rotated across the data rows (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, in column order: D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 45134; I = "Primera"; J = 50;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 3;  D = 44838; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 4;  D = 44838; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 5;  D = 45135; I = "Primera"; J = 70;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 6;  D = 44832; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 7;  D = 44832; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 8;  D = 44846; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 9;  D = 44846; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 10; D = 45133; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("P$row").Value = $r.P
}
